$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.252.51"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").Value = "1.879.10"
$ws.Range("E3").Value = "  -2.32%  "
$ws.Range("D5").Value = "'236.22"
$ws.Range("E5").Value = "  -1.57%  "
$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "'0.4843"
$ws.Range("E7").Value = "  -1.46%  "
$ws.Range("D8").Value = "'0.2875"
$ws.Range("E8").Value = "  -3.34%  "
$ws.Range("D9").Value = "'0.06590"
$ws.Range("E9").Value = "  -2.78%  "
$ws.Range("D10").Value = "1.880.99"
$ws.Range("E10").Value = "  -2.08%  "
$ws.Range("D11").Value = "'16.84"
$ws.Range("E11").Value = "  -1.23%  "
$ws.Range("D12").Value = "'0.07330"
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("D13").Value = "'5.117"
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("D14").Value = "'87.18"
$ws.Range("E14").Value = "  -3.29%  "
$ws.Range("D15").Value = "'0.6548"
$ws.Range("E15").Value = "  -3.06%  "
$ws.Range("D16").Value = "30.213.97"
$ws.Range("E16").Value = "  -2.02%  "
$ws.Range("E17").Value = "  -1.79%  "
$ws.Range("D18").Value = "'1.000"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").Value = "'0.000007746"
$ws.Range("E19").Value = "  -2.77%  "
$ws.Range("D20").Value = "'5.380"
$ws.Range("E20").Value = "  +3.93%  "
$ws.Range("D21").Value = "2.128.75"
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'193.13"
$ws.Range("E23").Value = "  -6.01%  "
$ws.Range("D24").Value = "'6.137"
$ws.Range("E24").Value = "  -2.06%  "
$ws.Range("D25").Value = "'9.260"
$ws.Range("E25").Value = "  -4.56%  "
$ws.Range("D26").Value = "'163.49"
$ws.Range("E26").Value = "  +2.94%  "
$ws.Range("D27").Value = "'18.04"
$ws.Range("E27").Value = "  -4.85%  "
$ws.Range("D28").Value = "'1.917"
$ws.Range("E28").Value = "  -3.79%  "
$ws.Range("D29").Value = "'1.434"
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("D30").Value = "'4.273"
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("D31").Value = "'0.09131"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("D32").Value = "'4.013"
$ws.Range("E32").Value = "  -1.48%  "
$ws.Range("D33").Value = "'0.05072"
$ws.Range("E33").Value = "  -2.56%  "
$ws.Range("D34").Value = "'0.7166"
$ws.Range("E34").Value = "  -5.12%  "
$ws.Range("D35").Value = "'1.103"
$ws.Range("E35").Value = "  -1.71%  "
$ws.Range("D36").Value = "'2.699"
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("D37").Value = "'0.01774"
$ws.Range("E37").Value = "  -5.01%  "
$ws.Range("D38").Value = "'2.638"
$ws.Range("E38").Value = "  -3.68%  "
$ws.Range("D39").Value = "'0.9237"
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("D40").Value = "'2.049"
$ws.Range("E40").Value = "  -2.14%  "
$ws.Range("D41").Value = "'105.74"
$ws.Range("E41").Value = "  -2.38%  "
$ws.Range("D42").Value = "'0.4275"
$ws.Range("E42").Value = "  -5.45%  "
$ws.Range("D43").Value = "'5.783"
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("D44").Value = "'0.9992"
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("D45").Value = "'7.411"
$ws.Range("E45").Value = "  -4.11%  "
$ws.Range("D46").Value = "'0.1305"
$ws.Range("E46").Value = "  -6.53%  "
$ws.Range("D47").Value = "'64.76"
$ws.Range("E47").Value = "  -8.44%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.501"
$ws.Range("E48").Value = "  +4.52%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.848"
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("D50").Value = "'33.82"
$ws.Range("E50").Value = "  -5.99%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05739"
$ws.Range("E51").Value = "  -3.55%  "
